# Add a "Save" column (H) to the s_vals sheet, mirroring the header style
# used by the existing columns (B..G) and filling in the per-row save flag.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: same text-header style as the existing headers (e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# Per-row "Save" values (0/1) for rows 2..51
$saveValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 1
    6  = 1
    7  = 1
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 1
    14 = 0
    15 = 1
    16 = 1
    17 = 1
    18 = 0
    19 = 0
    20 = 1
    21 = 0
    22 = 0
    23 = 0
    24 = 0
    25 = 1
    26 = 0
    27 = 1
    28 = 0
    29 = 1
    30 = 0
    31 = 0
    32 = 1
    33 = 0
    34 = 0
    35 = 0
    36 = 1
    37 = 1
    38 = 0
    39 = 0
    40 = 0
    41 = 0
    42 = 1
    43 = 0
    44 = 0
    45 = 0
    46 = 0
    47 = 0
    48 = 0
    49 = 0
    50 = 1
    51 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
